$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values in row 2 to reflect the latest FlashScore data refresh.
$ws.Range("G2").Value = 1.85
$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 2.63
$ws.Range("L2").Value = 5.5
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 2.5
$ws.Range("R2").Value = 1.5
$ws.Range("U2").Value = 2.25
$ws.Range("V2").Value = 1.57
$ws.Range("X2").Value = 7.5
$ws.Range("Z2").Value = 15
$ws.Range("AC2").Value = 7
$ws.Range("AH2").Value = 21
$ws.Range("AI2").Value = 17
$ws.Range("AJ2").Value = 51
$ws.Range("AN2").Value = 3.6
$ws.Range("AX2").Value = 29
